$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells retain their exact string representation
# (Excel would otherwise auto-convert numeric-looking strings to numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.688.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.58"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4836"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2892"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06550"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.828.86"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07465"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.119"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6686"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.648.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007590"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "233.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.108.37"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.282"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.188"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.395"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.65"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.76"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.959"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1025"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +11.80%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.347"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.036"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05077"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.212"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7514"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.18%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01885"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.648"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9216"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.067"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.91"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4298"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.637"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.425"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.29"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1280"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.491"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.038"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.47%  "
